$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column A, shifting columns B:F left to A:E
$ws.Columns("A").Delete()
